$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 37
$ws.Range("A37").Value = 112118082
$ws.Range("B37").Value = 89183
$ws.Range("E37").Value = 3215
$ws.Range("F37").Value = "Rödgul trumpetsvamp"
$ws.Range("G37").Value = "Craterellus lutescens"
$ws.Range("H37").Value = "(Fr.) Fr."
$ws.Range("I37").NumberFormat = "@"
$ws.Range("I37").Value = "25"
$ws.Range("I37").Style = "Normal"
$ws.Range("Q37").Value = 338493.1104975632
$ws.Range("R37").Value = 6433477.982833842
$ws.Range("AC37").Value = "1 mycel troligen"

# Row 38
$ws.Range("A38").Value = 112118106
$ws.Range("B38").Value = 89183
$ws.Range("E38").Value = 3215
$ws.Range("F38").Value = "Rödgul trumpetsvamp"
$ws.Range("G38").Value = "Craterellus lutescens"
$ws.Range("H38").Value = "(Fr.) Fr."
$ws.Range("I38").NumberFormat = "@"
$ws.Range("I38").Value = "10"
$ws.Range("I38").Style = "Normal"
$ws.Range("AC38").Value = "drygt 10, små."

# Row 39
$ws.Range("A39").Value = 112118071
$ws.Range("I39").NumberFormat = "@"
$ws.Range("I39").Value = "50"
$ws.Range("I39").Style = "Normal"
$ws.Range("Q39").Value = 338521.6570454396
$ws.Range("R39").Value = 6433487.434557982
$ws.Range("AC39").Value = "ca 50 f k- 2-3 mycel(?)"

# Row 40
$ws.Range("A40").Value = 112118115
$ws.Range("B40").Value = 95211
$ws.Range("E40").Value = 2606
$ws.Range("F40").Value = "Klippfrullania"
$ws.Range("G40").Value = "Frullania tamarisci"
$ws.Range("H40").Value = "(L.) Dumort."
$ws.Range("I40").ClearContents() | Out-Null
$ws.Range("J40").ClearContents() | Out-Null
$ws.Range("L40").ClearContents() | Out-Null
$ws.Range("P40").Value = "Angertuvan, söder om, Vg"
$ws.Range("Q40").Value = 338264.9079741923
$ws.Range("R40").Value = 6433236.348720711
$ws.Range("AC40").Value = "På klibbal i sumpskog."
$ws.Range("AJ40").Value = "klibbal"
$ws.Range("AK40").Value = "Alnus glutinosa"
$ws.Range("AO40").Value = "Alnus glutinosa"

# Row 42
$ws.Range("A42").Value = 112118117
$ws.Range("B42").Value = 89183
$ws.Range("E42").Value = 3215
$ws.Range("F42").Value = "Rödgul trumpetsvamp"
$ws.Range("G42").Value = "Craterellus lutescens"
$ws.Range("H42").Value = "(Fr.) Fr."
$ws.Range("L42").ClearContents() | Out-Null
$ws.Range("Q42").Value = 338295.1066844424
$ws.Range("R42").Value = 6433234.063291552
$ws.Range("AC42").Value = "En mindre fläck."
$ws.Range("AJ42").ClearContents() | Out-Null
$ws.Range("AK42").ClearContents() | Out-Null
$ws.Range("AO42").ClearContents() | Out-Null

# Row 43
$ws.Range("A43").Value = 112118111
$ws.Range("B43").Value = 89183
$ws.Range("E43").Value = 3215
$ws.Range("F43").Value = "Rödgul trumpetsvamp"
$ws.Range("G43").Value = "Craterellus lutescens"
$ws.Range("H43").Value = "(Fr.) Fr."
$ws.Range("I43").NumberFormat = "@"
$ws.Range("I43").Value = "20"
$ws.Range("I43").Style = "Normal"
$ws.Range("J43").Value = "fruktkroppar"
$ws.Range("L43").ClearContents() | Out-Null
$ws.Range("M43").ClearContents() | Out-Null
$ws.Range("Q43").Value = 338366.4565155458
$ws.Range("R43").Value = 6433512.816204711
$ws.Range("AC43").Value = "ca antal"

# Row 45
$ws.Range("A45").Value = 112118051
$ws.Range("B45").Value = 90332
$ws.Range("E45").Value = 4769
$ws.Range("F45").Value = "Svavelriska"
$ws.Range("G45").Value = "Lactarius scrobiculatus"
$ws.Range("H45").Value = "(Scop.:Fr.) Fr."
$ws.Range("I45").NumberFormat = "@"
$ws.Range("I45").Value = "7"
$ws.Range("I45").Style = "Normal"
$ws.Range("Q45").Value = 338499.1188164483
$ws.Range("R45").Value = 6433534.490804013
$ws.Range("AC45").Value = "3+4 f k några meter ifrån varandra - 2 mycel?"

# Row 46
$ws.Range("A46").Value = 112118103
$ws.Range("B46").Value = 90332
$ws.Range("E46").Value = 4769
$ws.Range("F46").Value = "Svavelriska"
$ws.Range("G46").Value = "Lactarius scrobiculatus"
$ws.Range("H46").Value = "(Scop.:Fr.) Fr."
$ws.Range("I46").NumberFormat = "@"
$ws.Range("I46").Value = "1"
$ws.Range("I46").Style = "Normal"
$ws.Range("J46").Value = "fruktkroppar"
$ws.Range("P46").Value = "Angertuvan, öster om, Vg"
$ws.Range("Q46").Value = 338374.6658049851
$ws.Range("R46").Value = 6433505.588431736
$ws.Range("AC46").ClearContents() | Out-Null

# Row 47
$ws.Range("A47").Value = 112118060
$ws.Range("B47").Value = 5135
$ws.Range("E47").Value = 105930
$ws.Range("F47").Value = "Vågbandad barkbock"
$ws.Range("G47").Value = "Semanotus undatus"
$ws.Range("H47").Value = "(Linnaeus, 1758)"
$ws.Range("I47").ClearContents() | Out-Null
$ws.Range("J47").ClearContents() | Out-Null
$ws.Range("L47").ClearContents() | Out-Null
$ws.Range("M47").Value = "äldre gnagspår"
$ws.Range("Q47").Value = 338513.2926211709
$ws.Range("R47").Value = 6433530.204112432
$ws.Range("AC47").Value = "I död gran."
